$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: merge "built in" (currently isolated between a gramStart/gramEnd
# proof-error pair) back into the surrounding plain-text run so it no longer
# sits in its own run with proofErr markers around it.
# ---------------------------------------------------------------------------

# Insert a sentinel character "Z" at the very start of the run that follows
# the "built in" proof-error pair (" view that allows ..."). This places the
# sentinel safely AFTER the closing gramEnd marker.
$rNext = $d.Content
$rNext.Find.Execute("view that allows", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$rNext.InsertBefore("Z")

# Delete "a built in Z" -- this span starts strictly before the gramStart
# marker and ends strictly after the gramEnd marker (because of the sentinel),
# so both now-interior proofErr markers are swept away along with the text.
$r2 = $d.Content
$r2.Find.Execute("a built in Z", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$r2.Delete()

# Re-insert the plain text that was removed, restoring the original wording
# but now as ordinary text with no enclosing proofErr run.
$r3 = $d.Range($r2.Start, $r2.Start)
$r3.InsertBefore("a built in ")

# ---------------------------------------------------------------------------
# Part 2: remove the trailing bold "Finish with anything" sentence (and its
# gramStart/gramEnd wrapper around "anything") that follows "...appointments. "
# ---------------------------------------------------------------------------

# Insert a sentinel character "Z" right after "anything" (and therefore after
# its closing gramEnd marker too, since "anything" is the last text in the
# paragraph and InsertAfter here creates a new trailing run).
$rTail = $d.Content
$rTail.Find.Execute("anything", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$rTail.InsertAfter("Z")

# Delete "Finish with anythingZ" -- starts strictly before the bold run and
# ends strictly after the gramEnd marker (because of the sentinel), so the
# bold text together with its proofErr wrapper is fully removed, leaving the
# paragraph ending cleanly after "...appointments. ".
$r4 = $d.Content
$r4.Find.Execute("Finish with anythingZ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$r4.Delete()
